$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 28375
$ws.Range("J125").Value = 36833.332
$ws.Range("L125").Value = 331499.988
$ws.Range("N125").Value = -336419.988

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2193.875
$ws.Range("I129").Value = 1573.3334
$ws.Range("J129").Value = 2991.7144
$ws.Range("K129").Value = 4720.0002
$ws.Range("L129").Value = 8975.143199999999
$ws.Range("M129").Value = 279.9997999999996
$ws.Range("N129").Value = -18975.1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2680.95
$ws.Range("I141").Value = 2680.95
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8042.849999999999
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -2862.849999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6333
$ws.Range("I45").Value = 5999
$ws.Range("K45").Value = 5999
$ws.Range("M45").Value = -5622

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2340.9473
$ws.Range("I74").Value = 2280.75
$ws.Range("K74").Value = 2280.75
$ws.Range("M74").Value = -1406.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2340.9473
$ws.Range("I77").Value = 2280.75
$ws.Range("K77").Value = 11403.75
$ws.Range("M77").Value = -7035.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 123178.6
$ws.Range("J135").Value = 123178.6
$ws.Range("L135").Value = 123178.6
$ws.Range("N135").Value = -133318.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 22335.25
$ws.Range("J88").Value = 22335.25
$ws.Range("L88").Value = 22335.25
$ws.Range("N88").Value = -23147.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 22335.25
$ws.Range("J91").Value = 22335.25
$ws.Range("L91").Value = 22335.25
$ws.Range("N91").Value = -25143.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1694.6136
$ws.Range("I134").Value = 1615.25
$ws.Range("K134").Value = 4845.75
$ws.Range("M134").Value = -2310.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 303737.25
$ws.Range("J137").Value = 238333
$ws.Range("L137").Value = 238333
$ws.Range("N137").Value = -248533

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 184.68182
$ws.Range("J33").Value = 249.92308
$ws.Range("L33").Value = 1499.53848
$ws.Range("N33").Value = -2065.53848

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 498.33334
$ws.Range("I39").Value = 497
$ws.Range("J39").Value = 499
$ws.Range("K39").Value = 1491
$ws.Range("L39").Value = 1497
$ws.Range("M39").Value = -1197
$ws.Range("N39").Value = -2085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1689.2222
$ws.Range("I68").Value = 300.5
$ws.Range("J68").Value = 2086
$ws.Range("K68").Value = 901.5
$ws.Range("L68").Value = 6258
$ws.Range("M68").Value = -90.5
$ws.Range("N68").Value = -7880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1689.2222
$ws.Range("I71").Value = 300.5
$ws.Range("J71").Value = 2086
$ws.Range("K71").Value = 2704.5
$ws.Range("L71").Value = 18774
$ws.Range("M71").Value = 1351.5
$ws.Range("N71").Value = -26886

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3400
$ws.Range("J80").Value = 3400
$ws.Range("L80").Value = 10200
$ws.Range("N80").Value = -12072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3400
$ws.Range("J83").Value = 3400
$ws.Range("L83").Value = 30600
$ws.Range("N83").Value = -39960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 100
$ws.Range("I86").Value = 150
$ws.Range("J86").Value = 75
$ws.Range("K86").Value = 450
$ws.Range("L86").Value = 225
$ws.Range("M86").Value = 736
$ws.Range("N86").Value = -2597

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 100
$ws.Range("I89").Value = 150
$ws.Range("J89").Value = 75
$ws.Range("K89").Value = 1350
$ws.Range("L89").Value = 675
$ws.Range("M89").Value = 4578
$ws.Range("N89").Value = -12531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 550
$ws.Range("I92").Value = 500
$ws.Range("K92").Value = 1500
$ws.Range("M92").Value = -252

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 1300
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2019.4615
$ws.Range("I107").Value = 3252.4285
$ws.Range("J107").Value = 581
$ws.Range("K107").Value = 9757.2855
$ws.Range("L107").Value = 1743
$ws.Range("M107").Value = -7837.2855
$ws.Range("N107").Value = -5583

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 3788.2856
$ws.Range("I112").Value = 2490.6
$ws.Range("J112").Value = 7032.5
$ws.Range("K112").Value = 7471.799999999999
$ws.Range("L112").Value = 21097.5
$ws.Range("M112").Value = -6363.799999999999
$ws.Range("N112").Value = -23313.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2826.5334
$ws.Range("I129").Value = 396.4
$ws.Range("K129").Value = 1189.2
$ws.Range("M129").Value = 3810.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 6082
$ws.Range("I130").Value = 6184
$ws.Range("K130").Value = 18552
$ws.Range("M130").Value = -13532

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1836.5555
$ws.Range("I131").Value = 712.2
$ws.Range("J131").Value = 2497.9412
$ws.Range("K131").Value = 2136.6
$ws.Range("L131").Value = 7493.823600000001
$ws.Range("M131").Value = 2903.4
$ws.Range("N131").Value = -17573.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3422.7856
$ws.Range("I134").Value = 1489.875
$ws.Range("K134").Value = 4469.625
$ws.Range("M134").Value = 600.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 4071
$ws.Range("J136").Value = 11500
$ws.Range("L136").Value = 34500
$ws.Range("N136").Value = -44700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6657
$ws.Range("J137").Value = 7867.1665
$ws.Range("L137").Value = 23601.4995
$ws.Range("N137").Value = -33801.49950000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3512.25
$ws.Range("J138").Value = 5000
$ws.Range("L138").Value = 15000
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 5116.6816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6469.55
$ws.Range("J80").Value = 6787.706
$ws.Range("L80").Value = 6787.706
$ws.Range("N80").Value = -8783.706

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6469.55
$ws.Range("J83").Value = 6787.706
$ws.Range("L83").Value = 33938.53
$ws.Range("N83").Value = -43922.53

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1568.8889
$ws.Range("I97").Value = 1662
$ws.Range("J97").Value = 1494.4
$ws.Range("K97").Value = 1662
$ws.Range("L97").Value = 1494.4
$ws.Range("M97").Value = -1166
$ws.Range("N97").Value = -2486.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 16874.5
$ws.Range("J136").Value = 16874.5
$ws.Range("L136").Value = 50623.5
$ws.Range("N136").Value = -55723.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3101.3333
$ws.Range("J46").Value = 3225.2942
$ws.Range("L46").Value = 3225.2942
$ws.Range("N46").Value = -3601.2942

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6784.2383
$ws.Range("I61").Value = 6437.75
$ws.Range("K61").Value = 6437.75
$ws.Range("M61").Value = -6235.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4919
$ws.Range("I68").Value = 2509.2222
$ws.Range("J68").Value = 6274.5
$ws.Range("K68").Value = 2509.2222
$ws.Range("L68").Value = 6274.5
$ws.Range("M68").Value = -1760.2222
$ws.Range("N68").Value = -7772.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4919
$ws.Range("I71").Value = 2509.2222
$ws.Range("J71").Value = 6274.5
$ws.Range("K71").Value = 12546.111
$ws.Range("L71").Value = 31372.5
$ws.Range("M71").Value = -8802.111000000001
$ws.Range("N71").Value = -38860.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6435.722
$ws.Range("I93").Value = 3977.5715
$ws.Range("K93").Value = 3977.5715
$ws.Range("M93").Value = -2729.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6784.2383
$ws.Range("I113").Value = 6437.75
$ws.Range("K113").Value = 6437.75
$ws.Range("M113").Value = -4267.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1599.8572
$ws.Range("I100").Value = 733.3333
$ws.Range("K100").Value = 1466.6666
$ws.Range("M100").Value = -925.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2249.5
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 6000
$ws.Range("M113").Value = -3830

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 173560
$ws.Range("J116").Value = 173560
$ws.Range("L116").Value = 173560
$ws.Range("N116").Value = -182738

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1107.3334
$ws.Range("I126").Value = 908.8333
$ws.Range("K126").Value = 2726.4999
$ws.Range("M126").Value = -256.4998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4326.155
$ws.Range("I132").Value = 4577.791
$ws.Range("J132").Value = 3334.4119
$ws.Range("K132").Value = 13733.373
$ws.Range("L132").Value = 10003.2357
$ws.Range("M132").Value = -11203.373
$ws.Range("N132").Value = -15063.2357
